# This workbook contains a single weekly data dump on Sheet1. A new week's
# observation is inserted at row 591 (directly under the header block),
# which pushes the existing data rows 591:690 down to 592:691 and grows
# the used range from A1:R690 to A1:R691.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 591; Excel shifts rows 591:690 down to 592:691 and
# carries the row-591 formatting (incl. the date style on column D) along.
$ws.Rows(591).Insert()

# Populate the newly inserted row with the new observation.
$ws.Range("A591").Value = 6
$ws.Range("B591").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C591").Value = "Metropolitana"
$ws.Range("D591").Value = 45218
$ws.Range("E591").Value = 13
$ws.Range("F591").Value = 100112043
$ws.Range("G591").Value = "Pepino ensalada"
$ws.Range("H591").Value = "Sin especificar"
$ws.Range("I591").Value = "Primera"
$ws.Range("J591").Value = 560
$ws.Range("K591").Value = 12000
$ws.Range("L591").Value = 13000
$ws.Range("M591").Value = 12429
$ws.Range("N591").Value = '$/caja 60 unidades'
$ws.Range("O591").Value = "Región de Arica y Parinacota"
$ws.Range("P591").Value = 207
$ws.Range("Q591").Value = 60
$ws.Range("R591").Value = "Hortaliza"
